# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (want-to-go count) figures and one
# "最低票价" (lowest price) cell that flipped to "不可售" (not for sale)
# on both the "展览" and "全部类型" sheets (which mirror the same rows).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        # row -> new F (想去人数) value
        $fUpdates = @{
            4  = 23
            6  = 15255
            7  = 412
            10 = 15307
            11 = 45
            12 = 8860
            13 = 358
            16 = 188
            20 = 38
            21 = 529
            33 = 34
            35 = 287
            37 = 111
            38 = 5429
            39 = 5228
        }
        $gRow = 2
    }
    else {
        # "全部类型" has 3 extra rows ahead of row 17, so later rows shift by +1
        $fUpdates = @{
            4  = 23
            6  = 15255
            7  = 412
            10 = 15307
            11 = 45
            12 = 8860
            13 = 358
            17 = 188
            21 = 38
            22 = 529
            36 = 34
            38 = 287
            40 = 111
            41 = 5429
            42 = 5228
        }
        $gRow = 2
    }

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # G2: 最低票价 20 -> "不可售" (no longer sellable)
    $ws.Cells.Item($gRow, 7).Value = "不可售"
}
